$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight PREMIUM = Yes rows with a yellow fill (added fill/style before
#     the text updates below, so the new style lands at the same index the
#     target workbook uses) ---
$ws.Range("E2").Interior.Color = 65535
$ws.Range("E3").Interior.Color = 65535
$ws.Range("E4").Interior.Color = 65535
$ws.Range("E7").Interior.Color = 65535

# --- Update data rows 2-16 with the latest scrape results ---
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1326041'
$ws.Range("C2").Value = 'ACE Program | Spanish Financial Analyst'
$ws.Range("D2").Value = 'Thane, Maharashtra, India'
$ws.Range("E2").Value = 'Yes'
$ws.Range("F2").Value = '30 applicants'
$ws.Range("G2").Value = '6 - 18 Months'
$ws.Range("H2").Value = 'Tata Consultancy Services Ltd.'

$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1308372'
$ws.Range("C3").Value = 'ACE Program | Talent Acquisition Specialist (Spanish Speaker)'
$ws.Range("D3").Value = 'Chennai, Tamil Nadu, India'
$ws.Range("E3").Value = 'Yes'
$ws.Range("F3").Value = '46 applicants'
$ws.Range("G3").Value = '6 - 18 Months'
$ws.Range("H3").Value = 'Tata Consultancy Services Ltd.'

$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1305153'
$ws.Range("C4").Value = 'ACE Program | Spanish Talent Acquisition Specialist'
$ws.Range("D4").Value = 'Chennai, Tamil Nadu, India'
$ws.Range("E4").Value = 'Yes'
$ws.Range("F4").Value = '55 applicants'
$ws.Range("G4").Value = '6 - 18 Months'
$ws.Range("H4").Value = 'Tata Consultancy Services Ltd.'

$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1328626'
$ws.Range("C5").Value = 'Travel Advisory intern'
$ws.Range("D5").Value = 'Hyderabad, Telangana, India'
$ws.Range("E5").Value = 'No'
$ws.Range("F5").Value = '1 applicant'
$ws.Range("G5").Value = '3 - 6 Months'
$ws.Range("H5").Value = 'Amaavi Luxe Travels'

$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1328565'
$ws.Range("C6").Value = 'WordPress Web Developer'
$ws.Range("D6").Value = 'Cairo, Cairo Governorate, Egypt'
$ws.Range("E6").Value = 'No'
$ws.Range("F6").Value = '0 applicants'
$ws.Range("G6").Value = '6 - 18 Months'
$ws.Range("H6").Value = 'AdMazad'

$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1327967'
$ws.Range("C7").Value = 'ACE Program | Russian Financial Analyst'
$ws.Range("D7").Value = 'Thane, Maharashtra, India'
$ws.Range("E7").Value = 'Yes'
$ws.Range("F7").Value = '5 applicants'
$ws.Range("G7").Value = '6 - 18 Months'
$ws.Range("H7").Value = 'Tata Consultancy Services Ltd.'

$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1327499'
$ws.Range("C8").Value = 'Full Stack Developer'
$ws.Range("D8").Value = 'Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt'
$ws.Range("E8").Value = 'No'
$ws.Range("F8").Value = '6 applicants'
$ws.Range("G8").Value = '3 - 6 Months'
$ws.Range("H8").Value = 'Skyline Egypt Tours'

$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1327497'
$ws.Range("C9").Value = 'SEO Specialist'
$ws.Range("D9").Value = 'Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt'
$ws.Range("E9").Value = 'No'
$ws.Range("F9").Value = '5 applicants'
$ws.Range("G9").Value = '3 - 6 Months'
$ws.Range("H9").Value = 'Skyline Egypt Tours'

$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1327432'
$ws.Range("C10").Value = 'Social Media Specialist'
$ws.Range("D10").Value = 'Zagazig, El-Hariry, Zagazig 1, Al-Sharqia Governorate, Egypt'
$ws.Range("E10").Value = 'No'
$ws.Range("F10").Value = '8 applicants'
$ws.Range("G10").Value = '6 - 18 Months'
$ws.Range("H10").Value = 'Admixy'

$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1327335'
$ws.Range("C11").Value = '2D Animator'
$ws.Range("D11").Value = 'Al Manteqah Ath Thamenah, Nasr City, Cairo Governorate, Egypt'
$ws.Range("E11").Value = 'No'
$ws.Range("F11").Value = '0 applicants'
$ws.Range("G11").Value = '9 - 12 Weeks'
$ws.Range("H11").Value = 'Sparks Studio'

$ws.Range("B12").Value = 'https://aiesec.org/opportunity/global-talent/1327221'
$ws.Range("C12").Value = 'Marketing Manager'
$ws.Range("D12").Value = 'Zagazig, El-Hariry, Zagazig 1, Al-Sharqia Governorate, Egypt'
$ws.Range("E12").Value = 'No'
$ws.Range("F12").Value = '14 applicants'
$ws.Range("G12").Value = '6 - 18 Months'
$ws.Range("H12").Value = 'Admixy'

$ws.Range("B13").Value = 'https://aiesec.org/opportunity/global-talent/1327220'
$ws.Range("C13").Value = 'Graphic designer'
$ws.Range("D13").Value = 'Zagazig, El-Hariry, Zagazig 1, Al-Sharqia Governorate, Egypt'
$ws.Range("E13").Value = 'No'
$ws.Range("F13").Value = '3 applicants'
$ws.Range("G13").Value = '6 - 18 Months'
$ws.Range("H13").Value = 'Admixy'

$ws.Range("B14").Value = 'https://aiesec.org/opportunity/global-talent/1327042'
$ws.Range("C14").Value = 'Content Creator'
$ws.Range("D14").Value = 'Sousse, Tunisia'
$ws.Range("E14").Value = 'No'
$ws.Range("F14").Value = '5 applicants'
$ws.Range("G14").Value = '6 - 18 Months'
$ws.Range("H14").Value = 'Progress Professional Center'

$ws.Range("B15").Value = 'https://aiesec.org/opportunity/global-talent/1326934'
$ws.Range("C15").Value = 'Social media manager'
$ws.Range("D15").Value = 'Alexandria, Alexandria Governorate, Egypt'
$ws.Range("E15").Value = 'No'
$ws.Range("F15").Value = '13 applicants'
$ws.Range("G15").Value = '6 - 18 Months'
$ws.Range("H15").Value = 'Eagle Office for Services'

$ws.Range("B16").Value = 'https://aiesec.org/opportunity/global-talent/1320868'
$ws.Range("C16").Value = 'Accelerate Romania|Data Labeling Specialist (SERBIAN Speackers)'
$ws.Range("D16").Value = 'Bucharest, Romania'
$ws.Range("E16").Value = 'No'
$ws.Range("F16").Value = '11 applicants'
$ws.Range("G16").Value = '9 - 12 Weeks'
$ws.Range("H16").Value = 'RepsMate'

# --- OPPORTUNITY ID column holds digit-only text; prefix with an apostrophe so
#     Excel keeps it as text instead of auto-converting it to a number, then
#     restore the Normal cell style so the stored style index is unaffected ---
$ws.Range("A2").Value = '''1326041'
$ws.Range("A3").Value = '''1308372'
$ws.Range("A4").Value = '''1305153'
$ws.Range("A5").Value = '''1328626'
$ws.Range("A6").Value = '''1328565'
$ws.Range("A7").Value = '''1327967'
$ws.Range("A8").Value = '''1327499'
$ws.Range("A9").Value = '''1327497'
$ws.Range("A10").Value = '''1327432'
$ws.Range("A11").Value = '''1327335'
$ws.Range("A12").Value = '''1327221'
$ws.Range("A13").Value = '''1327220'
$ws.Range("A14").Value = '''1327042'
$ws.Range("A15").Value = '''1326934'
$ws.Range("A16").Value = '''1320868'
$ws.Range("A2:A16").Style = "Normal"

# --- Remove the oldest listing (previously row 17); sheet now spans A1:H16 ---
$ws.Rows.Item(17).Delete()

# --- Widen columns C, D and narrow column H to fit the new content ---
$ws.Columns.Item(3).ColumnWidth = 65.16666666666667
$ws.Columns.Item(4).ColumnWidth = 63.166666666666664
$ws.Columns.Item(8).ColumnWidth = 32.166666666666664

